# Update gh-pages to output generated at 456a3b4
# Increments "想去人数" (F column) counts across sheets 展览, 演出, 全部类型
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5245
$ws.Range("F6").Value = 5245
$ws.Range("F7").Value = 151
$ws.Range("F9").Value = 530
$ws.Range("F11").Value = 1181
$ws.Range("F12").Value = 742
$ws.Range("F13").Value = 5175
$ws.Range("F15").Value = 72
$ws.Range("F16").Value = 90
$ws.Range("F17").Value = 288
$ws.Range("F18").Value = 288
$ws.Range("F19").Value = 254
$ws.Range("F21").Value = 251
$ws.Range("F22").Value = 3898
$ws.Range("F24").Value = 3817
$ws.Range("F25").Value = 185
$ws.Range("F28").Value = 232
$ws.Range("F29").Value = 246
$ws.Range("F32").Value = 112
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 6790
$ws.Range("F38").Value = 1101
$ws.Range("F39").Value = 523
$ws.Range("F40").Value = 102
$ws.Range("F42").Value = 59
$ws.Range("F43").Value = 1384
$ws.Range("F45").Value = 691
$ws.Range("F47").Value = 2315
$ws.Range("F50").Value = 782
$ws.Range("F51").Value = 929

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 7
$ws.Range("F7").Value = 134
$ws.Range("F9").Value = 89
$ws.Range("F16").Value = 7
$ws.Range("F25").Value = 814

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 5245
$ws.Range("F8").Value = 5245
$ws.Range("F9").Value = 151
$ws.Range("F11").Value = 89
$ws.Range("F12").Value = 530
$ws.Range("F13").Value = 1181
$ws.Range("F14").Value = 742
$ws.Range("F15").Value = 5175
$ws.Range("F17").Value = 72
$ws.Range("F18").Value = 90
$ws.Range("F19").Value = 288
$ws.Range("F20").Value = 288
$ws.Range("F21").Value = 254
$ws.Range("F23").Value = 251
$ws.Range("F24").Value = 3898
$ws.Range("F25").Value = 3817
$ws.Range("F26").Value = 185
$ws.Range("F28").Value = 232
$ws.Range("F29").Value = 246
$ws.Range("F32").Value = 112
$ws.Range("F35").Value = 18
$ws.Range("F37").Value = 6790
$ws.Range("F38").Value = 1101
$ws.Range("F39").Value = 523
$ws.Range("F41").Value = 102
$ws.Range("F43").Value = 59
$ws.Range("F44").Value = 1384
$ws.Range("F46").Value = 691
$ws.Range("F47").Value = 2315
$ws.Range("F49").Value = 782
$ws.Range("F50").Value = 929

